# Refresh the cryptos price/volume snapshot with the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.917.64"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.815.97"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.83"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4692"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07376"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8710"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "1.834.06"
$ws.Range("E12").Value = "  +3.74%  "
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.524"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07073"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.74"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008715"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.74"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "26.974.51"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.61"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").Value = "2.068.43"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.892"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.05"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.177"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.34"
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.342"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.28"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08959"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7682"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.510"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.910"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.089"
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05296"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.283"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5324"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.358"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1658"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.458"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4923"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.43"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.75"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06300"
$ws.Range("E51").Value = "  -0.48%  "
